$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$r = $t.Cell(1,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "14+2="
$r = $t.Cell(1,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "95-77="
$r = $t.Cell(1,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "9+18="
$r = $t.Cell(1,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "94-72="
$r = $t.Cell(1,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "17+27="
$r = $t.Cell(2,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "59-23="
$r = $t.Cell(2,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "36+14="
$r = $t.Cell(2,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "82-35="
$r = $t.Cell(2,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "84-54="
$r = $t.Cell(2,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "30+25="
$r = $t.Cell(3,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "52+1="
$r = $t.Cell(3,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "33+47="
$r = $t.Cell(3,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "16+55="
$r = $t.Cell(3,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "59-9="
$r = $t.Cell(3,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "54-27="
$r = $t.Cell(4,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "70-45="
$r = $t.Cell(4,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "62-0="
$r = $t.Cell(4,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "91-62="
$r = $t.Cell(4,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "56-49="
$r = $t.Cell(4,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "27-13="
$r = $t.Cell(5,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "33+48="
$r = $t.Cell(5,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "43-16="
$r = $t.Cell(5,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "43+7="
$r = $t.Cell(5,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "95-59="
$r = $t.Cell(5,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "12+65="
$r = $t.Cell(6,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "43+33="
$r = $t.Cell(6,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "89-59="
$r = $t.Cell(6,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "99-36="
$r = $t.Cell(6,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "41-10="
$r = $t.Cell(6,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "96-87="
$r = $t.Cell(7,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "81-26="
$r = $t.Cell(7,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "96-31="
$r = $t.Cell(7,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "17-14="
$r = $t.Cell(7,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49-30="
$r = $t.Cell(7,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "37+24="
$r = $t.Cell(8,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "94-87="
$r = $t.Cell(8,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "4+84="
$r = $t.Cell(8,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "29+50="
$r = $t.Cell(8,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "6+16="
$r = $t.Cell(8,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "75-5="
$r = $t.Cell(9,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "44-40="
$r = $t.Cell(9,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "74-49="
$r = $t.Cell(9,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "51-38="
$r = $t.Cell(9,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "8+54="
$r = $t.Cell(9,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "17+5="
$r = $t.Cell(10,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "45+1="
$r = $t.Cell(10,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "20+43="
$r = $t.Cell(10,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "24+47="
$r = $t.Cell(10,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "50-16="
$r = $t.Cell(10,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "82-73="
$r = $t.Cell(11,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "12+64="
$r = $t.Cell(11,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "17+61="
$r = $t.Cell(11,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "81-10="
$r = $t.Cell(11,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "46+26="
$r = $t.Cell(11,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "0+16="
$r = $t.Cell(12,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "51-44="
$r = $t.Cell(12,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "48+1="
$r = $t.Cell(12,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "47-13="
$r = $t.Cell(12,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "74-49="
$r = $t.Cell(12,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "10+6="
$r = $t.Cell(13,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "84+0="
$r = $t.Cell(13,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "4+71="
$r = $t.Cell(13,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "83-41="
$r = $t.Cell(13,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "24+37="
$r = $t.Cell(13,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "33+45="
$r = $t.Cell(14,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "45+34="
$r = $t.Cell(14,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "74-15="
$r = $t.Cell(14,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "56+38="
$r = $t.Cell(14,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "47+46="
$r = $t.Cell(14,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "93-75="
$r = $t.Cell(15,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "1+78="
$r = $t.Cell(15,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "13-9="
$r = $t.Cell(15,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "28-0="
$r = $t.Cell(15,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "79-47="
$r = $t.Cell(15,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "28+63="
$r = $t.Cell(16,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "67-18="
$r = $t.Cell(16,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "56-27="
$r = $t.Cell(16,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "48+29="
$r = $t.Cell(16,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "91-67="
$r = $t.Cell(16,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "74+14="
$r = $t.Cell(17,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "1+59="
$r = $t.Cell(17,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "53-43="
$r = $t.Cell(17,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "81-0="
$r = $t.Cell(17,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "68-15="
$r = $t.Cell(17,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "59+40="
$r = $t.Cell(18,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "19+67="
$r = $t.Cell(18,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "27+69="
$r = $t.Cell(18,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "98-97="
$r = $t.Cell(18,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "69+1="
$r = $t.Cell(18,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "32+16="
$r = $t.Cell(19,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "64+13="
$r = $t.Cell(19,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "54-6="
$r = $t.Cell(19,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "72-37="
$r = $t.Cell(19,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49-7="
$r = $t.Cell(19,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "47+14="
$r = $t.Cell(20,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "89+6="
$r = $t.Cell(20,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "57+14="
$r = $t.Cell(20,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "67-64="
$r = $t.Cell(20,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "25-3="
$r = $t.Cell(20,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "3+16="
